{"js": "const replacements = [\n  { oldText: \"766\u00d76=4596\", newText: \"549\u00d77=3843\" },\n  { oldText: \"977\u00d78=7816\", newText: \"539\u00d76=3234\" },\n  { oldText: \"948\u00d78=7584\", newText: \"834\u00d78=6672\" },\n  { oldText: \"330\u00d74=1320\", newText: \"174\u00d74=696\" },\n  { oldText: \"967\u00d75=4835\", newText: \"490\u00d73=1470\" },\n  { oldText: \"710\u00d77=4970\", newText: \"295\u00d75=1475\" },\n  { oldText: \"440\u00d73=1320\", newText: \"534\u00d73=1602\" },\n  { oldText: \"130\u00d76=780\", newText: \"260\u00d77=1820\" },\n  { oldText: \"520\u00d76=3120\", newText: \"383\u00d78=3064\" },\n  { oldText: \"360\u00d74=1440\", newText: \"346\u00d73=1038\" },\n  { oldText: \"445\u00d74=1780\", newText: \"934\u00d73=2802\" },\n  { oldText: \"543\u00d75=2715\", newText: \"472\u00d73=1416\" },\n  { oldText: \"722\u00d75=3610\", newText: \"536\u00d79=4824\" },\n  { oldText: \"618\u00d78=4944\", newText: \"211\u00d74=844\" },\n  { oldText: \"169\u00d73=507\", newText: \"195\u00d78=1560\" },\n  { oldText: \"821\u00d78=6568\", newText: \"268\u00d77=1876\" },\n  { oldText: \"939\u00d72=1878\", newText: \"259\u00d78=2072\" },\n  { oldText: \"871\u00d73=2613\", newText: \"682\u00d79=6138\" },\n  { oldText: \"212\u00d72=424\", newText: \"855\u00d78=6840\" },\n  { oldText: \"163\u00d74=652\", newText: \"642\u00d75=3210\" },\n  { oldText: \"873\u00d74=3492\", newText: \"309\u00d72=618\" },\n  { oldText: \"908\u00d75=4540\", newText: \"923\u00d74=3692\" },\n  { oldText: \"610\u00d75=3050\", newText: \"191\u00d73=573\" },\n  { oldText: \"805\u00d79=7245\", newText: \"701\u00d76=4206\" },\n  { oldText: \"798\u00d75=3990\", newText: \"135\u00d72=270\" },\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each \"Old\" value is the exact text of one three-digit x one-digit multiplication\n# expression (e.g. \"766x6=4596\") and occurs exactly once in the document, so a plain\n# exact-text Find/Replace is unambiguous for every pair below.\n$replacements = @(\n    @{ Old = \"766\u00d76=4596\"; New = \"549\u00d77=3843\" }\n    @{ Old = \"977\u00d78=7816\"; New = \"539\u00d76=3234\" }\n    @{ Old = \"948\u00d78=7584\"; New = \"834\u00d78=6672\" }\n    @{ Old = \"330\u00d74=1320\"; New = \"174\u00d74=696\" }\n    @{ Old = \"967\u00d75=4835\"; New = \"490\u00d73=1470\" }\n    @{ Old = \"710\u00d77=4970\"; New = \"295\u00d75=1475\" }\n    @{ Old = \"440\u00d73=1320\"; New = \"534\u00d73=1602\" }\n    @{ Old = \"130\u00d76=780\"; New = \"260\u00d77=1820\" }\n    @{ Old = \"520\u00d76=3120\"; New = \"383\u00d78=3064\" }\n    @{ Old = \"360\u00d74=1440\"; New = \"346\u00d73=1038\" }\n    @{ Old = \"445\u00d74=1780\"; New = \"934\u00d73=2802\" }\n    @{ Old = \"543\u00d75=2715\"; New = \"472\u00d73=1416\" }\n    @{ Old = \"722\u00d75=3610\"; New = \"536\u00d79=4824\" }\n    @{ Old = \"618\u00d78=4944\"; New = \"211\u00d74=844\" }\n    @{ Old = \"169\u00d73=507\"; New = \"195\u00d78=1560\" }\n    @{ Old = \"821\u00d78=6568\"; New = \"268\u00d77=1876\" }\n    @{ Old = \"939\u00d72=1878\"; New = \"259\u00d78=2072\" }\n    @{ Old = \"871\u00d73=2613\"; New = \"682\u00d79=6138\" }\n    @{ Old = \"212\u00d72=424\"; New = \"855\u00d78=6840\" }\n    @{ Old = \"163\u00d74=652\"; New = \"642\u00d75=3210\" }\n    @{ Old = \"873\u00d74=3492\"; New = \"309\u00d72=618\" }\n    @{ Old = \"908\u00d75=4540\"; New = \"923\u00d74=3692\" }\n    @{ Old = \"610\u00d75=3050\"; New = \"191\u00d73=573\" }\n    @{ Old = \"805\u00d79=7245\"; New = \"701\u00d76=4206\" }\n    @{ Old = \"798\u00d75=3990\"; New = \"135\u00d72=270\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    # MatchCase=True, Forward=True, Wrap=wdFindContinue(1), Replace=wdReplaceAll(2)\n    $find.Execute([ref]$pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$pair.New, 2) | Out-Null\n}\n"}
